$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper approach: for column D, values are numeric-looking price strings that
# must remain TEXT (matching source formatting, e.g. "1.002" is not the number 1.002).
# We force text number-format before writing, then restore the original cell style
# so no stray style/format difference is left behind.

# Row 2
$cell = $ws.Range('D2')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '29.765.12'
$cell.Style = $origStyle
$ws.Range('E2').Value = '  -1.36%  '

# Row 3
$cell = $ws.Range('D3')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.888.27'
$cell.Style = $origStyle
$ws.Range('E3').Value = '  -1.02%  '

# Row 4
$cell = $ws.Range('D4')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.002'
$cell.Style = $origStyle
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.7509'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  +2.60%  '

# Row 6
$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '238.93'
$cell.Style = $origStyle
$ws.Range('E6').Value = '  -1.97%  '

# Row 7
$cell = $ws.Range('D7')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = $origStyle
$ws.Range('E7').Value = '  -0.05%  '

# Row 8
$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.3030'
$cell.Style = $origStyle
$ws.Range('E8').Value = '  -3.31%  '

# Row 9
$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '25.37'
$cell.Style = $origStyle
$ws.Range('E9').Value = '  -5.37%  '

# Row 10
$cell = $ws.Range('D10')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.06789'
$cell.Style = $origStyle
$ws.Range('E10').Value = '  -1.76%  '

# Row 11
$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.07936'
$cell.Style = $origStyle
$ws.Range('E11').Value = '  -0.51%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.901.78'
$cell.Style = $origStyle
$ws.Range('E12').Value = '  -0.07%  '

# Row 13
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell = $ws.Range('D13')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.7402'
$cell.Style = $origStyle
$ws.Range('E13').Value = '  -4.98%  '

# Row 14
$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.148'
$cell.Style = $origStyle
$ws.Range('E14').Value = '  -2.01%  '

# Row 15
$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '90.39'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  -1.12%  '

# Row 16
$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '29.775.95'
$cell.Style = $origStyle
$ws.Range('E16').Value = '  -1.21%  '

# Row 17
$ws.Range('E17').Value = '  -2.62%  '

# Row 18
$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.898'
$cell.Style = $origStyle
$ws.Range('E18').Value = '  +0.70%  '

# Row 19
$cell = $ws.Range('D19')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '244.32'
$cell.Style = $origStyle
$ws.Range('E19').Value = '  +1.82%  '

# Row 20
$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.000007641'
$cell.Style = $origStyle
$ws.Range('E20').Value = '  -2.06%  '

# Row 21
$ws.Range('E21').Value = '  -0.01%  '

# Row 22
$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.003'
$cell.Style = $origStyle
$ws.Range('E22').Value = '  +0.07%  '

# Row 23
$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.869'
$cell.Style = $origStyle
$ws.Range('E23').Value = '  +0.77%  '

# Row 24
$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '9.192'
$cell.Style = $origStyle
$ws.Range('E24').Value = '  -2.17%  '

# Row 25
$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '165.16'
$cell.Style = $origStyle
$ws.Range('E25').Value = '  -0.18%  '

# Row 26
$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '18.61'
$cell.Style = $origStyle
$ws.Range('E26').Value = '  -2.67%  '

# Row 27
$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.1266'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  -0.29%  '

# Row 28
$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.005'
$cell.Style = $origStyle
$ws.Range('E28').Value = '  -3.82%  '

# Row 29
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.384'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  +2.42%  '

# Row 30
$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.509'
$cell.Style = $origStyle
$ws.Range('E30').Value = '  -2.43%  '

# Row 31
$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '4.227'
$cell.Style = $origStyle
$ws.Range('E31').Value = '  -1.74%  '

# Row 32
$cell = $ws.Range('D32')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '3.991'
$cell.Style = $origStyle
$ws.Range('E32').Value = '  -2.43%  '

# Row 33
$cell = $ws.Range('D33')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.05247'
$cell.Style = $origStyle
$ws.Range('E33').Value = '  +1.62%  '

# Row 34
$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.241'
$cell.Style = $origStyle
$ws.Range('E34').Value = '  -3.40%  '

# Row 35
$cell = $ws.Range('D35')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.7219'
$cell.Style = $origStyle
$ws.Range('E35').Value = '  -2.73%  '

# Row 36
$cell = $ws.Range('D36')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.715'
$cell.Style = $origStyle
$ws.Range('E36').Value = '  -1.30%  '

# Row 37
$cell = $ws.Range('D37')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.01897'
$cell.Style = $origStyle
$ws.Range('E37').Value = '  -2.47%  '

# Row 38
$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.772'
$cell.Style = $origStyle
$ws.Range('E38').Value = '  -0.97%  '

# Row 39
$cell = $ws.Range('D39')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.125'
$cell.Style = $origStyle
$ws.Range('E39').Value = '  -3.84%  '

# Row 40
$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.4373'
$cell.Style = $origStyle
$ws.Range('E40').Value = '  -1.57%  '

# Row 41
$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '71.54'
$cell.Style = $origStyle
$ws.Range('E41').Value = '  -4.12%  '

# Row 42
$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.001'
$cell.Style = $origStyle
$ws.Range('E42').Value = '  +0.07%  '

# Row 43
$cell = $ws.Range('D43')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.877'
$cell.Style = $origStyle
$ws.Range('E43').Value = '  -2.98%  '

# Row 44
$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.8240'
$cell.Style = $origStyle
$ws.Range('E44').Value = '  -1.26%  '

# Row 45
$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '100.04'
$cell.Style = $origStyle
$ws.Range('E45').Value = '  -1.21%  '

# Row 46
$cell = $ws.Range('D46')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.515'
$cell.Style = $origStyle
$ws.Range('E46').Value = '  -0.82%  '

# Row 47
$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '9.764'
$cell.Style = $origStyle
$ws.Range('E47').Value = '  +0.41%  '

# Row 48
$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.049.06'
$cell.Style = $origStyle
$ws.Range('E48').Value = '  +0.29%  '

# Row 49
$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '35.82'
$cell.Style = $origStyle
$ws.Range('E49').Value = '  -5.05%  '

# Row 50
$cell = $ws.Range('D50')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.05953'
$cell.Style = $origStyle
$ws.Range('E50').Value = '  -0.45%  '

# Row 51
$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.454'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  -0.25%  '
